$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows after row 336 first (inherits current format s=2/s=3 from row 336)
$ws.Rows("337:340").Insert()

# Fill in the new draw data
$data337 = @(2843,1,5,19,24,27,30,38,43,47,52,54,57,61,62,73,74,85,92,94,96)
$data338 = @(2844,1,2,5,14,16,26,31,32,34,49,51,54,61,66,70,77,91,95,96,97)
$data339 = @(2845,16,21,30,31,37,39,44,53,54,58,68,69,73,75,81,83,85,86,90,96)
$data340 = @(2846,12,13,14,16,25,28,33,39,41,42,53,60,67,68,73,74,76,92,93,97)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")

$row = 337
foreach ($data in @($data337,$data338,$data339,$data340)) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $data[$i]
    }
    $row++
}

# Now clear style of rows 326-336 (reverts residual purple/explicit style artifact)
$ws.Range("A326:U336").Style = "Normal"
